$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new date columns before the current newest-date column
# (old B, "Jun_17"), shifting old B:E ("Jun_17".."Jun_10") to E:H.
$ws.Range("B1:D1").EntireColumn.Insert()

# Keep the table's narrow 8-char column width (matching the pre-existing
# C:E formatting) across the whole, now-wider, C:H date range.
$ws.Range("C1:H1").EntireColumn.ColumnWidth = 7.1666666

# New header values for the two newest weekly pulls (B/C/D), the existing
# B:D header cells were pushed right to E:G by the insert above.
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# The newly inserted B:D columns need the same "UN" placeholder rating
# used across the rest of the table for rows 2-27.
$ws.Range("B2:D27").Value = "UN"

# Two new analyst rows appended at the bottom of the table.
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"
$ws.Range("D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"
$ws.Range("D29").Value = "UN"
